# Actualización automática 2025-09-01 08:30:07
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": zero out a handful of stray cell values ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("M6").Value2 = 0
$ws1.Range("M9").Value2 = 0
$ws1.Range("D10").Value2 = 0
$ws1.Range("L11").Value2 = 0
$ws1.Range("M12").Value2 = 0
$ws1.Range("M19").Value2 = 0
$ws1.Range("D20").Value2 = 0
$ws1.Range("D23").Value2 = 0
$ws1.Range("D24").Value2 = "0 de 22"
$ws1.Range("L24").Value2 = "0 de 22"
$ws1.Range("M24").Value2 = "0 de 22"

# --- Sheet "VENTA MENSUAL": shift the monthly columns left by one month ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Column widths follow the shift: old col E width moves to new col E position
# adjusted, and the new trailing month column widens.
# (ColumnWidth is stored with a constant +5/6 character padding offset, so we
# back that out here to land on the exact target OOXML <col width="..."/>.)
$ws2.Columns.Item(5).ColumnWidth = 13 - 0.8333333333333334
$ws2.Columns.Item(6).ColumnWidth = 16 - 0.8333333333333334

# Month headers shift left by one (mayo column is dropped, septiembre appended)
$ws2.Range("C1").Value2 = "junio"
$ws2.Range("D1").Value2 = "julio"
$ws2.Range("E1").Value2 = "agosto"
$ws2.Range("F1").Value2 = "septiembre"

for ($r = 2; $r -le 24; $r++) {
    $oldD = $ws2.Cells.Item($r, 4).Value2
    $oldE = $ws2.Cells.Item($r, 5).Value2
    $oldF = $ws2.Cells.Item($r, 6).Value2

    $ws2.Cells.Item($r, 3).Value2 = $oldD
    $ws2.Cells.Item($r, 4).Value2 = $oldE
    $ws2.Cells.Item($r, 5).Value2 = $oldF
    $ws2.Cells.Item($r, 6).Value2 = 0
}
